$d = $word.ActiveDocument
$rng = $d.Content
$rng.Find.Execute("Add random placement")
$paraStart = $rng.Paragraphs(1).Range.Start
$paraEnd = $rng.Paragraphs(1).Range.End
Write-Host "paraStart: " $paraStart " paraEnd: " $paraEnd
$fullRange = $d.Range($paraStart, $paraEnd)
Write-Host "fullRange start: " $fullRange.Start " end: " $fullRange.End
Write-Host "fullRange text: [" $fullRange.Text "]"
$fullRange.HighlightColorIndex = 4
Write-Host "Set highlight"
Write-Host "fullRange highlight after set: " $fullRange.HighlightColorIndex
